$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header cell H1 onto new headers I1, J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set new header labels
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Populate I and J columns for rows 2-60 with (I0, IF) pairs
$data = @{
    2 = @(7, 7)
    3 = @(5, 6)
    4 = @(9, 9)
    5 = @(4, 4)
    6 = @(8, 8)
    7 = @(9, 9)
    8 = @(4, 4)
    9 = @(9, 9)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(7, 7)
    13 = @(5, 5)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(4, 5)
    17 = @(7, 8)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(9, 10)
    23 = @(10, 10)
    24 = @(9, 9)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(6, 6)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(1, 1)
    32 = @(8, 9)
    33 = @(8, 8)
    34 = @(8, 8)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(7, 7)
    38 = @(5, 5)
    39 = @(8, 8)
    40 = @(7, 7)
    41 = @(9, 9)
    42 = @(6, 6)
    43 = @(5, 6)
    44 = @(8, 8)
    45 = @(8, 8)
    46 = @(3, 3)
    47 = @(5, 5)
    48 = @(5, 5)
    49 = @(8, 8)
    50 = @(6, 6)
    51 = @(6, 6)
    52 = @(8, 8)
    53 = @(4, 4)
    54 = @(9, 9)
    55 = @(8, 8)
    56 = @(6, 6)
    57 = @(5, 5)
    58 = @(8, 8)
    59 = @(7, 7)
    60 = @(8, 8)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value2 = $vals[0]
    $ws.Cells.Item($r, 10).Value2 = $vals[1]
}

